# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" worksheet right after "总计" (i.e. before the
#    existing "2022-Q3" sheet), populated with that quarter's fund-holding
#    breakdown.
# 2) Insert a new top data row in "总计" for 2022-Q4, pushing the older
#    quarters down by one row.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# --- 1. New "2022-Q4" sheet, inserted before "2022-Q3" -------------------
$newSheet = $wb.Worksheets.Add($q3Sheet, $null)
$newSheet.Name = "2022-Q4"

# Clone header-row (s=2) and index-column (s=2) look-and-feel from the
# existing "2022-Q3" sheet so the new sheet's styling matches exactly.
$q3Sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q3Sheet.Range("A2:A11").Copy($newSheet.Range("A2:A14"))

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$rows = @(
    @(0, "519674", "银河创新成长混合A", "145.89", "92.48", "7.94", "11.5837", 7),
    @(1, "014143", "银河创新成长混合C", "22.12", "92.48", "7.94", "1.7563", 7),
    @(2, "016105", "申万菱信兴乐优选混合A", "3.37", "54.93", "3.77", "0.1270", 6),
    @(3, "016106", "申万菱信兴乐优选混合C", "2.90", "54.93", "3.77", "0.1093", 6),
    @(4, "013340", "创金合信芯片产业股票C", "0.94", "92.41", "7.46", "0.0701", 1),
    @(5, "000522", "华润元大信息传媒科技混合", "1.38", "62.01", "4.98", "0.0687", 3),
    @(6, "013339", "创金合信芯片产业股票A", "0.92", "92.41", "7.46", "0.0686", 1),
    @(7, "012200", "新华鑫科技3个月滚动持有灵活配置混合A", "1.13", "93.39", "3.64", "0.0411", 6),
    @(8, "004890", "中邮健康文娱灵活配置混合", "0.42", "92.60", "4.62", "0.0194", 5),
    @(9, "002772", "光大保德信产业新动力灵活配置混合", "0.27", "92.57", "6.48", "0.0175", 6),
    @(10, "012201", "新华鑫科技3个月滚动持有灵活配置混合C", "0.42", "93.39", "3.64", "0.0153", 6),
    @(11, "000531", "东吴阿尔法灵活配置混合A", "0.26", "84.60", "3.45", "0.0090", 10),
    @(12, "014581", "东吴阿尔法灵活配置混合C", "0.03", "84.60", "3.45", "0.0010", 10)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# --- 2. New first data row on "总计" for 2022-Q4 --------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 13
$totalSheet.Cells.Item(2, 4).Value = 13.89
